$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing C-column counts (no date/age change) ---
$ws.Cells.Item(972, 3).Value = 10
$ws.Cells.Item(979, 3).Value = 37
$ws.Cells.Item(988, 3).Value = 2
$ws.Cells.Item(1007, 3).Value = 30
$ws.Cells.Item(1011, 3).Value = 37
$ws.Cells.Item(1020, 3).Value = 36
$ws.Cells.Item(1024, 3).Value = 17
$ws.Cells.Item(1030, 3).Value = 37
$ws.Cells.Item(1036, 3).Value = 37
$ws.Cells.Item(1040, 3).Value = 38
$ws.Cells.Item(1043, 3).Value = 12
$ws.Cells.Item(1045, 3).Value = 25

# --- Rows 1057-1065 get re-mapped age groups / dates / counts ---
$ws.Cells.Item(1057, 2).Value = "50-59"
$ws.Cells.Item(1057, 3).Value = 2

$ws.Cells.Item(1058, 2).Value = "60-69"
$ws.Cells.Item(1058, 3).Value = 13

$ws.Cells.Item(1059, 2).Value = "70-79"
$ws.Cells.Item(1059, 3).Value = 18

$ws.Cells.Item(1060, 2).Value = "80+"
$ws.Cells.Item(1060, 3).Value = 39

$ws.Cells.Item(1061, 1).Value = 44175
$ws.Cells.Item(1061, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1061, 2).Value = "40-49"
$ws.Cells.Item(1061, 3).Value = 1

$ws.Cells.Item(1062, 3).Value = 4
$ws.Cells.Item(1063, 3).Value = 11
$ws.Cells.Item(1064, 3).Value = 18
$ws.Cells.Item(1065, 3).Value = 27

# --- Append new rows 1066-1074 ---
$newRows = @(
    @(44176, "30-39", 1),
    @(44176, "40-49", 2),
    @(44176, "50-59", 3),
    @(44176, "60-69", 5),
    @(44176, "70-79", 17),
    @(44176, "80+", 24),
    @(44177, "60-69", 1),
    @(44177, "70-79", 3),
    @(44177, "80+", 8)
)

$r = 1066
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
